# Update the "想去人数" (interested count) numbers in the F column across
# the sheets of the workbook, matching the upstream data refresh recorded
# in the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1160
$ws.Range("F3").Value  = 1973
$ws.Range("F4").Value  = 624
$ws.Range("F5").Value  = 1278
$ws.Range("F9").Value  = 345
$ws.Range("F10").Value = 129
$ws.Range("F11").Value = 104
$ws.Range("F12").Value = 864
$ws.Range("F13").Value = 264
$ws.Range("F14").Value = 137
$ws.Range("F19").Value = 709
$ws.Range("F20").Value = 84
$ws.Range("F21").Value = 674
$ws.Range("F22").Value = 207
$ws.Range("F24").Value = 918
$ws.Range("F25").Value = 376
$ws.Range("F26").Value = 201
$ws.Range("F28").Value = 310

# --- 演出 (Performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 337

# --- 本地生活 (Local life) sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 332

# --- 全部类型 (All types, combined) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 332
$ws.Range("F3").Value  = 1160
$ws.Range("F4").Value  = 1973
$ws.Range("F5").Value  = 624
$ws.Range("F6").Value  = 1278
$ws.Range("F11").Value = 345
$ws.Range("F12").Value = 129
$ws.Range("F13").Value = 104
$ws.Range("F14").Value = 864
$ws.Range("F15").Value = 264
$ws.Range("F16").Value = 137
$ws.Range("F19").Value = 337
$ws.Range("F26").Value = 709
$ws.Range("F27").Value = 84
$ws.Range("F28").Value = 674
$ws.Range("F29").Value = 207
$ws.Range("F31").Value = 918
$ws.Range("F32").Value = 376
$ws.Range("F35").Value = 201
$ws.Range("F37").Value = 310

$wb.Save()
